$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1423738854105352
$ws.Range("D2").Value = 0.1356897205440646
$ws.Range("E2").Value = 0.151826455498675
$ws.Range("F2").Value = 1.919499299024324
$ws.Range("G2").Value = 1.265468727497094
$ws.Range("H2").Value = 1.184471871255852
$ws.Range("J2").Value = 0.1994485699222608
$ws.Range("K2").Value = 1.124391946476948
$ws.Range("N2").Value = 1.707544937395062
$ws.Range("B3").Value = 0.1329632149027731
$ws.Range("D3").Value = 0.1325294824970484
$ws.Range("E3").Value = 0.1484496169491365
$ws.Range("F3").Value = 1.909255574359747
$ws.Range("G3").Value = 1.256115846258211
$ws.Range("H3").Value = 1.185978054678841
$ws.Range("J3").Value = 0.1951129908603804
$ws.Range("K3").Value = 1.023407134790062
$ws.Range("N3").Value = 1.728115892529334
$ws.Range("B4").Value = 0.127255543272085
$ws.Range("D4").Value = 0.1306409127408301
$ws.Range("E4").Value = 0.1464543052680902
$ws.Range("F4").Value = 1.904108472138518
$ws.Range("G4").Value = 1.251236117740234
$ws.Range("H4").Value = 1.187515970769397
$ws.Range("J4").Value = 0.1925691845663735
$ws.Range("K4").Value = 0.9617477853782646
$ws.Range("N4").Value = 1.741381535198564
$ws.Range("B5").Value = 0.1249474987936736
$ws.Range("D5").Value = 0.1298843885376471
$ws.Range("E5").Value = 0.1456608338512915
$ws.Range("F5").Value = 1.902297863419875
$ws.Range("G5").Value = 1.249464019701819
$ws.Range("H5").Value = 1.188296725971767
$ws.Range("J5").Value = 0.1915622447554099
$ws.Range("K5").Value = 0.9367082429581899
$ws.Range("N5").Value = 1.746946992958385
$ws.Range("B6").Value = 0.1245653338490911
$ws.Range("D6").Value = 0.1297595599648815
$ws.Range("E6").Value = 0.1455302644152425
$ws.Range("F6").Value = 1.902014527297794
$ws.Range("G6").Value = 1.24918281942827
$ws.Range("H6").Value = 1.188435669778229
$ws.Range("J6").Value = 0.1913968343034966
$ws.Range("K6").Value = 0.9325557156461457
$ws.Range("N6").Value = 1.747880770379781
$ws.Range("B7").Value = 0.1272243435960974
$ws.Range("D7").Value = 0.1306306569491369
$ws.Range("E7").Value = 0.1464435247351759
$ws.Range("F7").Value = 1.904082892663837
$ws.Range("G7").Value = 1.251211343046776
$ws.Range("H7").Value = 1.18752587677649
$ws.Range("J7").Value = 0.1925554845135622
$ws.Range("K7").Value = 0.9614097402337904
$ws.Range("N7").Value = 1.741455946793423
$ws.Range("B8").Value = 0.1391145343527569
$ws.Range("D8").Value = 0.1345893422090541
$ws.Range("E8").Value = 0.1506459144505214
$ws.Range("F8").Value = 1.915729847062906
$ws.Range("G8").Value = 1.262064391081481
$ws.Range("H8").Value = 1.184863869710782
$ws.Range("J8").Value = 0.1979290830905427
$ws.Range("K8").Value = 1.089500714410491
$ws.Range("N8").Value = 1.714505967381674
$ws.Range("B9").Value = 0.162985595212831
$ws.Range("D9").Value = 0.1427618679804254
$ws.Range("E9").Value = 0.1595069565112937
$ws.Range("F9").Value = 1.947659097121701
$ws.Range("G9").Value = 1.290223412604632
$ws.Range("H9").Value = 1.184515730084371
$ws.Range("J9").Value = 0.2094083070025192
$ws.Range("K9").Value = 1.3434349002261
$ws.Range("N9").Value = 1.666698012180141
$ws.Range("B10").Value = 0.1808567350223029
$ws.Range("D10").Value = 0.1490144185243736
$ws.Range("E10").Value = 0.1663970390670713
$ws.Range("F10").Value = 1.976697465245365
$ws.Range("G10").Value = 1.315148206913051
$ws.Range("H10").Value = 1.187242151782357
$ws.Range("J10").Value = 0.2184218920943977
$ws.Range("K10").Value = 1.531703051563397
$ws.Range("N10").Value = 1.634648814541492
$ws.Range("B11").Value = 0.1890581734868562
$ws.Range("D11").Value = 0.1519125016754685
$ws.Range("E11").Value = 0.1696144561886328
$ws.Range("F11").Value = 1.991128201505902
$ws.Range("G11").Value = 1.327416941880983
$ws.Range("H11").Value = 1.189132815806516
$ws.Range("J11").Value = 0.2226496302155283
$ws.Range("K11").Value = 1.617728110633038
$ws.Range("N11").Value = 1.620737093653511
$ws.Range("B12").Value = 0.1921740431186976
$ws.Range("D12").Value = 0.1530176243370249
$ws.Range("E12").Value = 0.1708447741622763
$ws.Range("F12").Value = 1.996768943000149
$ws.Range("G12").Value = 1.332197297956156
$ws.Range("H12").Value = 1.189942487777273
$ws.Range("J12").Value = 0.2242689705043546
$ws.Range("K12").Value = 1.650358498295134
$ws.Range("N12").Value = 1.615565184351899
$ws.Range("B13").Value = 0.1915025353500681
$ws.Range("D13").Value = 0.1527792755560853
$ws.Range("E13").Value = 0.1705792715296468
$ws.Range("F13").Value = 1.995546265899108
$ws.Range("G13").Value = 1.331161772157571
$ws.Range("H13").Value = 1.189763939184246
$ws.Range("J13").Value = 0.2239193977490856
$ws.Range("K13").Value = 1.643328536967829
$ws.Range("N13").Value = 1.616674768416205
$ws.Range("B14").Value = 0.1893143151785068
$ws.Range("D14").Value = 0.1520032671185305
$ws.Range("E14").Value = 0.1697154356260455
$ws.Range("F14").Value = 1.991588735518505
$ws.Range("G14").Value = 1.327807526150679
$ws.Range("H14").Value = 1.189197548638589
$ws.Range("J14").Value = 0.2227824854204385
$ws.Range("K14").Value = 1.620411539885311
$ws.Range("N14").Value = 1.620309669522644
$ws.Range("B15").Value = 0.1879752859797463
$ws.Range("D15").Value = 0.1515289384202987
$ws.Range("E15").Value = 0.1691878678320649
$ws.Range("F15").Value = 1.989187588316369
$ws.Range("G15").Value = 1.325770483364096
$ws.Range("H15").Value = 1.188862828739786
$ws.Range("J15").Value = 0.222088490369174
$ws.Range("K15").Value = 1.606381317271314
$ws.Range("N15").Value = 1.622548678612129
$ws.Range("B16").Value = 0.1803221824395678
$ws.Range("D16").Value = 0.148826099601763
$ws.Range("E16").Value = 0.1661884455313754
$ws.Range("F16").Value = 1.975779000051475
$ws.Range("G16").Value = 1.31436519411281
$ws.Range("H16").Value = 1.187131696479611
$ws.Range("J16").Value = 0.2181481689449782
$ws.Range("K16").Value = 1.526088758534968
$ws.Range("N16").Value = 1.635571425542799
$ws.Range("B17").Value = 0.1756455208371648
$ws.Range("D17").Value = 0.1471817306381809
$ws.Range("E17").Value = 0.164369679978293
$ws.Range("F17").Value = 1.967866392508128
$ws.Range("G17").Value = 1.307607173631197
$ws.Range("H17").Value = 1.186236421542333
$ws.Range("J17").Value = 0.2157635894470076
$ws.Range("K17").Value = 1.476929285531469
$ws.Range("N17").Value = 1.643731560953127
$ws.Range("B18").Value = 0.172962393260164
$ws.Range("D18").Value = 0.1462409968316791
$ws.Range("E18").Value = 0.1633313940934471
$ws.Range("F18").Value = 1.963430166324002
$ws.Range("G18").Value = 1.303807664620336
$ws.Range("H18").Value = 1.185782690212733
$ws.Range("J18").Value = 0.2144040311288791
$ws.Range("K18").Value = 1.448689924121481
$ws.Range("N18").Value = 1.648487884997164
$ws.Range("B19").Value = 0.1720550989585803
$ws.Range("D19").Value = 0.1459233517625194
$ws.Range("E19").Value = 0.1629811912874786
$ws.Range("F19").Value = 1.961947853249939
$ws.Range("G19").Value = 1.302536228337914
$ws.Range("H19").Value = 1.185639571190706
$ws.Range("J19").Value = 0.2139457651291394
$ws.Range("K19").Value = 1.439134728464865
$ws.Range("N19").Value = 1.650109081079307
$ws.Range("B20").Value = 0.1761426604878977
$ws.Range("D20").Value = 0.1473562528103685
$ws.Range("E20").Value = 0.1645624813096163
$ws.Range("F20").Value = 1.968696807762171
$ws.Range("G20").Value = 1.308317513006955
$ws.Range("H20").Value = 1.186325388993538
$ws.Range("J20").Value = 0.216016190743531
$ws.Range("K20").Value = 1.482158684808326
$ws.Range("N20").Value = 1.642856396827974
$ws.Range("B21").Value = 0.1899567739224324
$ws.Range("D21").Value = 0.1522309915322921
$ws.Range("E21").Value = 0.1699688407594806
$ws.Range("F21").Value = 1.992746373447872
$ws.Range("G21").Value = 1.328789095118964
$ws.Range("H21").Value = 1.189361366232163
$ws.Range("J21").Value = 0.2231159245673382
$ws.Range("K21").Value = 1.627141338214301
$ws.Range("N21").Value = 1.619239400015084
$ws.Range("B22").Value = 0.19904425963189
$ws.Range("D22").Value = 0.1554616656601411
$ws.Range("E22").Value = 0.1735718767357071
$ws.Range("F22").Value = 2.009491064106058
$ws.Range("G22").Value = 1.342952465172147
$ws.Range("H22").Value = 1.191891907602582
$ws.Range("J22").Value = 0.2278632192894037
$ws.Range("K22").Value = 1.72221383771415
$ws.Range("N22").Value = 1.60436488169038
$ws.Range("B23").Value = 0.1941887359326557
$ws.Range("D23").Value = 0.1537333158218672
$ws.Range("E23").Value = 0.1716424923898998
$ws.Range("F23").Value = 2.000459961491984
$ws.Range("G23").Value = 1.33532124743607
$ws.Range("H23").Value = 1.190491254354413
$ws.Range("J23").Value = 0.2253196672293996
$ws.Range("K23").Value = 1.671442830107367
$ws.Range("N23").Value = 1.612252349048864
$ws.Range("B24").Value = 0.1759178863603523
$ws.Range("D24").Value = 0.1472773368936515
$ws.Range("E24").Value = 0.1644752929455677
$ws.Range("F24").Value = 1.968321025572436
$ws.Range("G24").Value = 1.307996101458372
$ws.Range("H24").Value = 1.186284976878738
$ws.Range("J24").Value = 0.2159019542938267
$ws.Range("K24").Value = 1.479794401508968
$ws.Range("N24").Value = 1.643251856201522
$ws.Range("B25").Value = 0.1564689939619512
$ws.Range("D25").Value = 0.1405072922466957
$ws.Range("E25").Value = 0.1570432320749049
$ws.Range("F25").Value = 1.938044309231884
$ws.Range("G25").Value = 1.281864929796399
$ws.Range("H25").Value = 1.184087089772703
$ws.Range("J25").Value = 0.2062014735515305
$ws.Range("K25").Value = 1.274441976903006
$ws.Range("N25").Value = 1.679091455372128
